# Updated cryptos list on Tue Jan 16 05:12:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "price" cells hold plain numeric-looking text (e.g. "315.55") in the
# original workbook (t="inlineStr"). A bare .Value assignment lets Excel
# auto-coerce those strings into real numbers, which would flip the cell type.
# Forcing NumberFormat to Text ("@") first keeps the assigned value stored as text,
# matching the source data which is a formatted string, not a numeric type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.797.50"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.529.64"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.55"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.87"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.55"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.917.89"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.528.84"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.853.12"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.94"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.43"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.68"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.29"
$ws.Range("E29").Value = "  +5.94%  "
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.53"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.19"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.34"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0781"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.28"
$ws.Range("E40").Value = "  -3.96%  "
$ws.Range("E41").Value = "  +13.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.030.19"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.18"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.01"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.71"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.87"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.772.07"
$ws.Range("E51").Value = "  +0.49%  "
